$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.901.38'
$ws.Range('E2').Value = '  +0.11%  '
$ws.Range('D3').Value = '1.892.17'
$ws.Range('E3').Value = '  -0.11%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7719'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.61%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '243.64'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.10%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.001'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3124'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.69%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '25.65'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.18%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07220'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.04%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08700'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +7.43%  '
$ws.Range('D12').Value = '2.024.40'
$ws.Range('E12').Value = '  +6.04%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.7703'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.61%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.402'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.57%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '94.22'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +1.91%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.209'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.66%  '
$ws.Range('D17').Value = '30.044.91'
$ws.Range('E17').Value = '  +0.58%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.91'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.41%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '245.27'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.55%  '
$ws.Range('B20').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C20').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D20').Value = '2.292.63'
$ws.Range('E20').Value = '  +6.01%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.000007861'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.91%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '8.171'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.34%  '
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1592'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -3.46%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.515'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.88%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '162.64'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.18%  '
$ws.Range('E28').Value = '  +0.35%  '
$ws.Range('E29').Value = '  -0.52%  '
$ws.Range('E30').Value = '  +1.68%  '
$ws.Range('E31').Value = '  -0.24%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.514'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.21%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.117'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.25%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05457'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.89%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.248'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.75%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7524'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.20%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.003'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.44%  '
$ws.Range('E38').Value = '  +2.85%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01986'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +3.18%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.783'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.01%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4510'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.93%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '73.93'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.22%  '
$ws.Range('D43').Value = '1.100.36'
$ws.Range('E43').Value = '  -3.83%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.077'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +3.62%  '
$ws.Range('E45').Value = '  +0.61%  '
$ws.Range('B46').Value = 'PaxDollar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.000'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.04%  '
$ws.Range('B47').Value = 'RocketPoolETH'
$ws.Range('C47').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D47').Value = '2.177.30'
$ws.Range('E47').Value = '  +5.92%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '102.95'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.05%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.883'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.05%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.616'
$ws.Range('D50').ClearFormats()
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '9.850'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.46%  '
